$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") rows 2-28 changed from serial date 45435 (2024-05-23)
# to serial date 45437 (2024-05-25).
for ($row = 2; $row -le 28; $row++) {
    $ws.Cells.Item($row, 3).Value = 45437
}
